# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (fund-holding detail, same shape as the
# existing quarter sheets) right after "总计" and before "2021-Q4", and adds
# a matching summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "2022-Q3" sheet, right before the "2021-Q4" sheet so the
#    tab order becomes: 总计, 2022-Q3, 2021-Q4, 2021-Q2
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($q4)
$newSheet.Name = "2022-Q3"

# Header row (bold, centered, bordered - matches the other quarter sheets)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$headerRange = $newSheet.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$aRange = $newSheet.Range("A2:A3")
$aRange.Font.Bold = $true
$aRange.HorizontalAlignment = -4108
$aRange.VerticalAlignment = -4160
$aRange.Borders.LineStyle = 1
$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1

# Row 2 (B,D,E,F,G hold text-like values - e.g. fund code "011351" has a
# leading zero, so force a text number-format before writing the value so
# it isn't silently coerced into a number).
$textRow2 = $newSheet.Range("B2:B2,D2:G2")
$textRow2.NumberFormat = "@"
$newSheet.Range("B2").Value = "011351"
$newSheet.Range("C2").Value = "金鹰年年邮益一年持有期混合A"
$newSheet.Range("D2").Value = "3.43"
$newSheet.Range("E2").Value = "34.33"
$newSheet.Range("F2").Value = "0.87"
$newSheet.Range("G2").Value = "0.0298"
$newSheet.Range("H2").Value = 6

# Row 3
$textRow3 = $newSheet.Range("B3:B3,D3:G3")
$textRow3.NumberFormat = "@"
$newSheet.Range("B3").Value = "011352"
$newSheet.Range("C3").Value = "金鹰年年邮益一年持有期混合C"
$newSheet.Range("D3").Value = "0.27"
$newSheet.Range("E3").Value = "34.33"
$newSheet.Range("F3").Value = "0.87"
$newSheet.Range("G3").Value = "0.0023"
$newSheet.Range("H3").Value = 6

# ---------------------------------------------------------------------
# 2. Insert the matching summary row into "总计" (right after the header,
#    pushing the existing 2021-Q4 / 2021-Q2 rows down by one).
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.03

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2

Write-Host "2022-Q3 sheet added"
